$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("G2").Value = 0.398441
$ws.Range("H2").Value = 0.796882
$ws.Range("I2").Value = 0.08945363909080989
$ws.Range("J2").Value = 0.06146862341190577
$ws.Range("K2").Value = 2
$ws.Range("M2").Value = 17.945525
$ws.Range("N2").Value = 35.89105
$ws.Range("O2").Value = 0.3147738875783
$ws.Range("P2").Value = 0.2472168478181395
$ws.Range("Q2").Value = 7.150232926525
$ws.Range("R2").Value = 28.6009317061
$ws.Range("S2").Value = 0.02815766973464042
$ws.Range("T2").Value = 0.01519607931961164

$ws.Range("E3").Value = 2
$ws.Range("G3").Value = 0.398441
$ws.Range("H3").Value = 0.796882
$ws.Range("I3").Value = 0.08945363909080989
$ws.Range("J3").Value = 0.06146862341190577
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.914977
$ws.Range("N3").Value = 11.744931
$ws.Range("O3").Value = 0.06867074270993077
$ws.Range("P3").Value = 0.0808988541617353
$ws.Range("Q3").Value = 1.559887350857
$ws.Range("R3").Value = 9.359324105142
$ws.Range("S3").Value = 0.006142847834472011
$ws.Range("T3").Value = 0.004972741200922393

$ws.Range("E4").Value = 2
$ws.Range("G4").Value = 0.398441
$ws.Range("H4").Value = 0.796882
$ws.Range("I4").Value = 0.08945363909080989
$ws.Range("J4").Value = 0.06146862341190577
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 7.595080333333333
$ws.Range("N4").Value = 22.785241
$ws.Range("O4").Value = 0.133221678551774
$ws.Range("P4").Value = 0.1569442927079769
$ws.Range("Q4").Value = 3.026191403093666
$ws.Range("R4").Value = 18.157148418562
$ws.Range("S4").Value = 0.01191716395224228
$ws.Range("T4").Value = 0.00964714962511454

$ws.Range("E5").Value = 2
$ws.Range("G5").Value = 0.398441
$ws.Range("H5").Value = 0.796882
$ws.Range("I5").Value = 0.08945363909080989
$ws.Range("J5").Value = 0.06146862341190577
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 8.835736333333331
$ws.Range("N5").Value = 26.507209
$ws.Range("O5").Value = 0.1549834332102386
$ws.Range("P5").Value = 0.18258113522554
$ws.Range("Q5").Value = 3.520519620389666
$ws.Range("R5").Value = 21.123117722338
$ws.Range("S5").Value = 0.01386383209944332
$ws.Range("T5").Value = 0.01122301104329696

$ws.Range("E6").Value = 2
$ws.Range("G6").Value = 0.398441
$ws.Range("H6").Value = 0.796882
$ws.Range("I6").Value = 0.08945363909080989
$ws.Range("J6").Value = 0.06146862341190577
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 10.81295533333333
$ws.Range("N6").Value = 32.438866
$ws.Range("O6").Value = 0.1896648878471846
$ws.Range("P6").Value = 0.2234382721964117
$ws.Range("Q6").Value = 4.308324735968667
$ws.Range("R6").Value = 25.849948415812
$ws.Range("S6").Value = 0.01696621442568098
$ws.Range("T6").Value = 0.01373444300944813

$ws.Range("E7").Value = 2
$ws.Range("G7").Value = 0.398441
$ws.Range("H7").Value = 0.796882
$ws.Range("I7").Value = 0.08945363909080989
$ws.Range("J7").Value = 0.06146862341190577
$ws.Range("K7").Value = 2
$ws.Range("M7").Value = 7.90657
$ws.Range("N7").Value = 15.81314
$ws.Range("O7").Value = 0.1386853701025721
$ws.Range("P7").Value = 0.1089205978901965
$ws.Range("Q7").Value = 3.15030165737
$ws.Range("R7").Value = 12.60120662948
$ws.Range("S7").Value = 0.01240591104433088
$ws.Range("T7").Value = 0.006695199213512104

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 4.055721
$ws.Range("H8").Value = 12.167163
$ws.Range("I8").Value = 0.9105463609091901
$ws.Range("J8").Value = 0.9385313765880943
$ws.Range("K8").Value = 2
$ws.Range("M8").Value = 17.945525
$ws.Range("N8").Value = 35.89105
$ws.Range("O8").Value = 0.3147738875783
$ws.Range("P8").Value = 0.2472168478181395
$ws.Range("Q8").Value = 72.782042598525
$ws.Range("R8").Value = 436.69225559115
$ws.Range("S8").Value = 0.2866162178436596
$ws.Range("T8").Value = 0.2320207684985279

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 4.055721
$ws.Range("H9").Value = 12.167163
$ws.Range("I9").Value = 0.9105463609091901
$ws.Range("J9").Value = 0.9385313765880943
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 3.914977
$ws.Range("N9").Value = 11.744931
$ws.Range("O9").Value = 0.06867074270993077
$ws.Range("P9").Value = 0.0808988541617353
$ws.Range("Q9").Value = 15.878054433417
$ws.Range("R9").Value = 142.902489900753
$ws.Range("S9").Value = 0.06252789487545876
$ws.Range("T9").Value = 0.0759261129608129

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 4.055721
$ws.Range("H10").Value = 12.167163
$ws.Range("I10").Value = 0.9105463609091901
$ws.Range("J10").Value = 0.9385313765880943
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 7.595080333333333
$ws.Range("N10").Value = 22.785241
$ws.Range("O10").Value = 0.133221678551774
$ws.Range("P10").Value = 0.1569442927079769
$ws.Range("Q10").Value = 30.803526804587
$ws.Range("R10").Value = 277.231741241283
$ws.Range("S10").Value = 0.1213045145995317
$ws.Range("T10").Value = 0.1472971430828623

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 4.055721
$ws.Range("H11").Value = 12.167163
$ws.Range("I11").Value = 0.9105463609091901
$ws.Range("J11").Value = 0.9385313765880943
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 8.835736333333331
$ws.Range("N11").Value = 26.507209
$ws.Range("O11").Value = 0.1549834332102386
$ws.Range("P11").Value = 0.18258113522554
$ws.Range("Q11").Value = 35.83528139756299
$ws.Range("R11").Value = 322.5175325780669
$ws.Range("S11").Value = 0.1411196011107952
$ws.Range("T11").Value = 0.1713581241822431

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 4.055721
$ws.Range("H12").Value = 12.167163
$ws.Range("I12").Value = 0.9105463609091901
$ws.Range("J12").Value = 0.9385313765880943
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 10.81295533333333
$ws.Range("N12").Value = 32.438866
$ws.Range("O12").Value = 0.1896648878471846
$ws.Range("P12").Value = 0.2234382721964117
$ws.Range("Q12").Value = 43.854330017462
$ws.Range("R12").Value = 394.6889701571581
$ws.Range("S12").Value = 0.1726986734215036
$ws.Range("T12").Value = 0.2097038291869636

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 4.055721
$ws.Range("H13").Value = 12.167163
$ws.Range("I13").Value = 0.9105463609091901
$ws.Range("J13").Value = 0.9385313765880943
$ws.Range("K13").Value = 2
$ws.Range("M13").Value = 7.90657
$ws.Range("N13").Value = 15.81314
$ws.Range("O13").Value = 0.1386853701025721
$ws.Range("P13").Value = 0.1089205978901965
$ws.Range("Q13").Value = 32.06684198697
$ws.Range("R13").Value = 192.40105192182
$ws.Range("S13").Value = 0.1262794590582412
$ws.Range("T13").Value = 0.1022253986766844
